{"js": "// Append new paragraphs to the end of the document body, after the last\n// existing paragraph (\" g:out out, outside, out of\"), mirroring the\n// target diff: a blank line, several single-run lines, another blank\n// line, then a block of \"# \"-prefixed two-run comment/code lines.\n\nconst body = context.document.body;\nconst paragraphs = body.paragraphs;\nparagraphs.load(\"items\");\nawait context.sync();\n\nlet anchor = paragraphs.items[paragraphs.items.length - 1];\n\n// Simple single-run paragraphs (in order).\nconst simpleLines = [\n  \"\",\n  \"# test ability to add inflections after the fact\",\n  \"\\\\inflect n akin = number PCL\",\n  \"ja = n boat\",\n  \"\\\\inflect v eme = negation MAYBE\",\n  \"taril = v jump\",\n  \"\",\n];\n\nfor (const line of simpleLines) {\n  anchor = anchor.insertParagraph(line, Word.InsertLocation.after);\n  await context.sync();\n}\n\n// Paragraphs made of two runs: a \"# \" lead-in run followed by the rest\n// of the line as a second run (both carry the same non-italic\n// formatting, but are kept as distinct runs).\nconst twoRunLines = [\n  [\"# \", \"# should throw error if you try to inflect or instantiate an undefined POS\"],\n  [\"# \", \"\\\\inflect nx morphx = featurex valx\"],\n  [\"# \", \"xxx = nx thing\"],\n  [\"# \", \"\\\\pos nx _{featurex}\"],\n  [\"# \", \"# the morph \\\"morphx\\\" and the lexeme \\\"xxx\\\" should not have been added even though the POS \\\"nx\\\" was eventually defined\"],\n  [\"# \", \"\\\\inflect nx morph2x = featurex val2x\"],\n  [\"# \", \"x2x = nx thing2\"],\n];\n\nfor (const [first, second] of twoRunLines) {\n  anchor = anchor.insertParagraph(first, Word.InsertLocation.after);\n  await context.sync();\n\n  const tail = anchor.getRange(Word.RangeLocation.end);\n  const secondRun = tail.insertText(second, Word.InsertLocation.end);\n  // Briefly toggle italic on/off so the engine keeps this insertion as\n  // its own run instead of silently merging it into the previous run\n  // (both runs end up with identical, explicit non-italic formatting).\n  secondRun.font.set({ italic: true });\n  await context.sync();\n  secondRun.font.set({ italic: false });\n  await context.sync();\n}\n", "ps1": "$d = $word.ActiveDocument\n\n# Simple single-run paragraphs appended (in order) after the current\n# last paragraph (\" g:out out, outside, out of\").\n$lines = @(\n  \"\",\n  \"# test ability to add inflections after the fact\",\n  \"\\inflect n akin = number PCL\",\n  \"ja = n boat\",\n  \"\\inflect v eme = negation MAYBE\",\n  \"taril = v jump\",\n  \"\"\n)\n\nforeach ($line in $lines) {\n  $r = $d.Paragraphs.Last.Range\n  $r.InsertParagraphAfter()\n  $p = $d.Paragraphs.Last\n  if ($line -ne \"\") {\n    $p.Range.Text = $line\n  }\n}\n\n# Paragraphs made of two runs: a \"# \" lead-in run followed by the rest\n# of the line as a second run (both end up with identical, explicit\n# non-italic formatting, but stay as distinct <w:r> runs).\n$pairs = @(\n  @(\"# \", \"# should throw error if you try to inflect or instantiate an undefined POS\"),\n  @(\"# \", \"\\inflect nx morphx = featurex valx\"),\n  @(\"# \", \"xxx = nx thing\"),\n  @(\"# \", \"\\pos nx _{featurex}\"),\n  @(\"# \", \"# the morph `\"morphx`\" and the lexeme `\"xxx`\" should not have been added even though the POS `\"nx`\" was eventually defined\"),\n  @(\"# \", \"\\inflect nx morph2x = featurex val2x\"),\n  @(\"# \", \"x2x = nx thing2\")\n)\n\nforeach ($pair in $pairs) {\n  $first = $pair[0]\n  $second = $pair[1]\n\n  $r = $d.Paragraphs.Last.Range\n  $r.InsertParagraphAfter()\n  $p = $d.Paragraphs.Last\n  $p.Range.Text = $first\n\n  $pEnd = $p.Range.End\n  $tail = $d.Range($pEnd, $pEnd)\n  $tail.InsertAfter($second)\n  # Briefly toggle italic on/off so the engine keeps this insertion as\n  # its own run instead of silently merging it into the previous run.\n  $tail.Font.Italic = 1\n  $tail.Font.Italic = 0\n}\n"}
